$wb = $excel.ActiveWorkbook

# --- OFF sheet (sheet1) ---
$wsOff = $wb.Worksheets.Item("OFF")

$wsOff.Range("B2").Value = 228
$wsOff.Range("C2").Value = 170
$wsOff.Range("D2").Value = 59
$wsOff.Range("E2").Value = 26

$wsOff.Range("B3").Value = 260
$wsOff.Range("C3").Value = 194
$wsOff.Range("D3").Value = 69
$wsOff.Range("E3").Value = 37
$wsOff.Range("F3").Value = 3

# --- DEF sheet (sheet2) ---
$wsDef = $wb.Worksheets.Item("DEF")

$wsDef.Range("B2").Value = 274
$wsDef.Range("C2").Value = 192
$wsDef.Range("D2").Value = 55
$wsDef.Range("E2").Value = 26
$wsDef.Range("F2").Value = 3

$wsDef.Range("B3").Value = 209
$wsDef.Range("C3").Value = 142
$wsDef.Range("D3").Value = 42
$wsDef.Range("E3").Value = 22
